$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Shared-string rich-text edits (report title / date-range strings)
# ---------------------------------------------------------------------------
# "Volume 30   Number  30" -> "Volume 30   Number  31"  (cell A8)
$ws.Range("A8").Characters(21, 2).Text = "31"

# "Report Covering the Week  7/24/2023  Through  7/30/2023" -> "... 7/31/2023 ... 8/6/2023" (cell C9)
$ws.Range("C9").Characters(27, 9).Text = "7/31/2023"
$ws.Range("C9").Characters(47, 9).Text = "8/6/2023"

# ---------------------------------------------------------------------------
# Row 15
# ---------------------------------------------------------------------------
$ws.Range("D15").Value = 5
$ws.Range("G15").Value = 7
$ws.Range("J15").Value = 15
$ws.Range("K15").Value = 20

# ---------------------------------------------------------------------------
# Row 16  (D16/E16 flip from text placeholders to real numbers)
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("D15").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = 2
$ws.Range("E15").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 40
$ws.Range("I16").Value = 156
$ws.Range("J16").Value = 113
$ws.Range("K16").Value = 38.053097345132
$ws.Range("L16").Value = 126.086956521739
$ws.Range("M16").Value = 47.169811320754
$ws.Range("N16").Value = -75.663026521060

# ---------------------------------------------------------------------------
# Row 17  (D17/E17 flip from text placeholders to real numbers)
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 5
$ws.Range("D15").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = 5
$ws.Range("E15").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 50
$ws.Range("I17").Value = 143
$ws.Range("J17").Value = 135
$ws.Range("K17").Value = 5.925925925925
$ws.Range("L17").Value = 26.548672566371
$ws.Range("M17").Value = 88.157894736842
$ws.Range("N17").Value = -22.282608695652

# ---------------------------------------------------------------------------
# Row 18
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 21
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 75
$ws.Range("I18").Value = 117
$ws.Range("J18").Value = 103
$ws.Range("K18").Value = 13.592233009708
$ws.Range("L18").Value = 17
$ws.Range("M18").Value = -20.408163265306
$ws.Range("N18").Value = -86.473988439306

# ---------------------------------------------------------------------------
# Row 19
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -38.461538461538
$ws.Range("F19").Value = 58
$ws.Range("G19").Value = 56
$ws.Range("H19").Value = 3.571428571428
$ws.Range("I19").Value = 442
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 10.5
$ws.Range("L19").Value = 88.085106382978
$ws.Range("M19").Value = 64.312267657992
$ws.Range("N19").Value = -15.809523809523

# ---------------------------------------------------------------------------
# Row 20
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 125
$ws.Range("F20").Value = 27
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = 28.571428571428
$ws.Range("I20").Value = 147
$ws.Range("J20").Value = 114
$ws.Range("K20").Value = 28.947368421052
$ws.Range("L20").Value = 31.25
$ws.Range("M20").Value = 9.701492537313
$ws.Range("N20").Value = -87.542372881355

# ---------------------------------------------------------------------------
# Row 21
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -18.75
$ws.Range("G21").Value = 122
$ws.Range("H21").Value = 18.032786885245
$ws.Range("I21").Value = 1023
$ws.Range("J21").Value = 880
$ws.Range("K21").Value = 16.25
$ws.Range("L21").Value = 60.849056603773
$ws.Range("M21").Value = 37.5
$ws.Range("N21").Value = -69.982394366197

# ---------------------------------------------------------------------------
# Row 22  (D22/E22 flip from text placeholders to real numbers)
# ---------------------------------------------------------------------------
$ws.Range("C22").Value = 1
$ws.Range("D15").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$ws.Range("E15").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = 0
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 600
$ws.Range("I22").Value = 55
$ws.Range("J22").Value = 44
$ws.Range("K22").Value = 25
$ws.Range("L22").Value = 175
$ws.Range("M22").Value = 120

# ---------------------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 54
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 80
$ws.Range("F24").Value = 206
$ws.Range("G24").Value = 124
$ws.Range("H24").Value = 66.129032258064
$ws.Range("I24").Value = 1155
$ws.Range("J24").Value = 813
$ws.Range("K24").Value = 42.066420664206
$ws.Range("L24").Value = 62.905500705218
$ws.Range("M24").Value = 105.516014234875

# ---------------------------------------------------------------------------
# Row 25
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 12.5
$ws.Range("F25").Value = 31
$ws.Range("G25").Value = 49
$ws.Range("H25").Value = -36.734693877551
$ws.Range("I25").Value = 304
$ws.Range("J25").Value = 350
$ws.Range("K25").Value = -13.142857142857
$ws.Range("L25").Value = 12.592592592592
$ws.Range("M25").Value = 0

# ---------------------------------------------------------------------------
# Row 26
# ---------------------------------------------------------------------------
$ws.Range("D26").Value = 5
$ws.Range("G26").Value = 8
$ws.Range("J26").Value = 20
$ws.Range("K26").Value = 15
$ws.Range("L26").Value = 76.923076923076

# ---------------------------------------------------------------------------
# Row 27  (D27/E27 flip from text placeholders to real numbers)
# ---------------------------------------------------------------------------
$ws.Range("C27").Value = 1
$ws.Range("C27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 3
$ws.Range("H27").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -12.5
$ws.Range("I27").Value = 56
$ws.Range("J27").Value = 56
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 55.555555555555

# ---------------------------------------------------------------------------
# Row 30  (D30/E30 flip from real numbers to text placeholders)
# ---------------------------------------------------------------------------
$ws.Range("D30").Value = "'0"
$ws.Range("C30").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = "***.*"
$ws.Range("C30").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("L30").Value = -50
